$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.434.96'
$ws.Range('E2').Value = '  -2.20%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.986.29'
$ws.Range('E3').Value = '  -6.22%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '329.94'
$ws.Range('E5').Value = '  -4.79%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.007'
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4964'
$ws.Range('E7').Value = '  -4.61%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4194'
$ws.Range('E8').Value = '  -6.17%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '51.92'
$ws.Range('E9').Value = '  -3.76%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.08848'
$ws.Range('E10').Value = '  -5.47%  '
$ws.Range('E11').Value = '  -5.56%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '23.28'
$ws.Range('E12').Value = '  -8.23%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '8.020'
$ws.Range('E13').Value = '  -7.45%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.487'
$ws.Range('E14').Value = '  -7.02%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '1.930.25'
$ws.Range('E15').Value = '  -7.77%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '95.98'
$ws.Range('E16').Value = '  -6.58%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.007'
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('E18').Value = '  -5.56%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06627'
$ws.Range('E19').Value = '  -1.23%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '19.68'
$ws.Range('E20').Value = '  -8.79%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.005'
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.944'
$ws.Range('E22').Value = '  -5.82%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '29.436.13'
$ws.Range('E23').Value = '  -2.27%  '
$ws.Range('E24').Value = '  -7.15%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.282'
$ws.Range('E25').Value = '  -2.18%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '157.54'
$ws.Range('E26').Value = '  -3.59%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.52'
$ws.Range('E27').Value = '  -7.50%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.505'
$ws.Range('E28').Value = '  -3.70%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.329'
$ws.Range('E29').Value = '  -8.72%  '
$ws.Range('E30').Value = '  -4.83%  '
$ws.Range('E31').Value = '  -9.34%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.09911'
$ws.Range('E32').Value = '  -6.44%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.565'
$ws.Range('E33').Value = '  -12.36%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.840'
$ws.Range('E34').Value = '  -7.13%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.788'
$ws.Range('E35').Value = '  -4.63%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '9.535'
$ws.Range('E36').Value = '  -11.66%  '
$ws.Range('E37').Value = '  -7.83%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06339'
$ws.Range('E38').Value = '  -7.90%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.285'
$ws.Range('E39').Value = '  -3.44%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.6495'
$ws.Range('E40').Value = '  -9.03%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '11.72'
$ws.Range('E41').Value = '  -7.98%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.2063'
$ws.Range('E42').Value = '  -8.36%  '
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.6317'
$ws.Range('E44').Value = '  -9.58%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.209'
$ws.Range('E45').Value = '  -8.02%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '13.41'
$ws.Range('E46').Value = '  -9.09%  '
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.531'
$ws.Range('E48').Value = '  -2.80%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.06981'
$ws.Range('E49').Value = '  -2.99%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.140'
$ws.Range('E50').Value = '  -6.28%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.00000000322'
$ws.Range('E51').Value = '  -7.82%  '
